$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.221.62'
$ws.Range("E2").Value = '  +0.56%  '
$ws.Range("D3").Value = '3.130.63'
$ws.Range("E3").Value = '  +0.87%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.56%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.127.22'
$ws.Range("E8").Value = '  +0.83%  '
$ws.Range("E9").Value = '  +0.57%  '
$ws.Range("E10").Value = '  +0.07%  '
$ws.Range("E11").Value = '  +2.96%  '
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("E13").Value = '  +2.63%  '
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").Value = '3.643.68'
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("E16").Value = '  +2.91%  '
$ws.Range("D17").Value = '64.204.58'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").Value = '3.098.29'
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '479.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.58'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.709'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.68'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.67'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.77%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  -0.79%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.34'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.15'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.28%  '
$ws.Range("E30").Value = '  -3.87%  '
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.81'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.53%  '
$ws.Range("E34").Value = '  -2.86%  '
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.97'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.82%  '
$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D37").Value = '0.0₃0770'
$ws.Range("E37").Value = '  +3.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.38'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.86%  '
$ws.Range("E39").Value = '  +2.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '444.69'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.55%  '
$ws.Range("E41").Value = '  +0.32%  '
$ws.Range("E42").Value = '  +0.73%  '
$ws.Range("E43").Value = '  -1.59%  '
$ws.Range("D44").Value = '2.854.04'
$ws.Range("E44").Value = '  +0.59%  '
$ws.Range("E45").Value = '  -2.19%  '
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("E47").Value = '  +1.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.54%  '
$ws.Range("E50").Value = '  +0.34%  '
